$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '98.112.97'
$ws.Range('E2').Value = '  +4.08%  '
$ws.Range('D3').Value = '3.353.03'
$ws.Range('E3').Value = '  +9.11%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '253.94'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +7.12%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '620.49'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.66%  '
$ws.Range('E7').Value = '  +7.32%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.383'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.46%  '
$ws.Range('E9').Value = '  -0.02%  '
$ws.Range('D10').Value = '3.348.84'
$ws.Range('E10').Value = '  +9.03%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.790'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.44%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.198'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.29%  '
$ws.Range('D13').Value = '97.832.77'
$ws.Range('E13').Value = '  +4.05%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.77'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +6.11%  '
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '3.982.51'
$ws.Range('E15').Value = '  +9.31%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000245'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.76%  '
$ws.Range('E17').Value = '  +2.59%  '
$ws.Range('D18').Value = '3.352.73'
$ws.Range('E18').Value = '  +9.60%  '
$ws.Range('E19').Value = '  +1.16%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.77'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '479.44'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +9.32%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.82'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.85%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.0000206'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +8.65%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.09'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.03%  '
$ws.Range('E25').Value = '  +2.18%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '87.52'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.31%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.91'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.62%  '
$ws.Range('D28').Value = '3.571.15'
$ws.Range('E28').Value = '  +10.35%  '
$ws.Range('E29').Value = '  -0.18%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.187'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +5.07%  '
$ws.Range('E31').Value = '  -0.94%  '
$ws.Range('E32').Value = '  -0.65%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.00'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.95%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '9.18'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.54%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '27.15'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +6.96%  '
$ws.Range('B36').Value = 'Bittensor'
$ws.Range('C36').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '517.19'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +7.30%  '
$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.150'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.12%  '
$ws.Range('E38').Value = '  -7.15%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.94'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.78%  '
$ws.Range('E40').Value = '  +3.24%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.446'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.31%  '
$ws.Range('E42').Value = '  -1.45%  '
$ws.Range('E43').Value = '  -3.85%  '
$ws.Range('B44').Value = 'ARBITRUM'
$ws.Range('C44').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.788'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +17.32%  '
$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.22'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.29%  '
$ws.Range('E46').Value = '  -0.01%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '160.71'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.46%  '
$ws.Range('E48').Value = '  +5.46%  '
$ws.Range('B49').Value = 'ImmutableX'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.37'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +6.83%  '
$ws.Range('B50').Value = 'OKB'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '45.44'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.15%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.47'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +5.37%  '
